# work on SWA and local sleep section
#
# 1. The "Body" paragraph switches from the "First Paragraph" style to the
#    "Body Text" style.
# 2. Several paragraph styles (Heading 1/2/3, Body Text, Bibliography) swap
#    their Georgia-based fonts for the "CMU Serif" family.

$d = $word.ActiveDocument

# --- 1. Re-style the "Body" paragraph -------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Style.NameLocal -eq "First Paragraph") {
        $p.Style = "Body Text"
    }
}

# --- 2. Update fonts used by the template's heading / body styles ---------

# Heading 1: Georgia -> CMU Serif (ascii/hAnsi/cs)
$h1 = $d.Styles("Heading 1")
$h1.Font.Name = "CMU Serif"
$h1.Font.NameBi = "CMU Serif"

# Heading 2: Georgia -> CMU Serif (ascii/hAnsi/cs)
$h2 = $d.Styles("Heading 2")
$h2.Font.Name = "CMU Serif"
$h2.Font.NameBi = "CMU Serif"

# Heading 3: Georgia -> CMU SERIF BOLDITALIC (ascii/hAnsi/cs)
$h3 = $d.Styles("Heading 3")
$h3.Font.Name = "CMU SERIF BOLDITALIC"
$h3.Font.NameBi = "CMU SERIF BOLDITALIC"

# Body Text: Georgia / Times New Roman (Body CS) -> CMU Serif Roman
$bt = $d.Styles("Body Text")
$bt.Font.Name = "CMU Serif Roman"
$bt.Font.NameBi = "CMU Serif Roman"

# Bibliography: add CMU Serif Roman (ascii/hAnsi/cs)
$bib = $d.Styles("Bibliography")
$bib.Font.Name = "CMU Serif Roman"
$bib.Font.NameBi = "CMU Serif Roman"
